$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text number format on price cells whose new values would
# otherwise be auto-parsed as numeric literals by Excel (single "." as
# decimal separator), so they remain text matching the source data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update price (D) and volume-change (E) cells per latest crypto data pull.
$ws.Range("D2").Value = "27.684.16"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.872.54"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "331.86"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "0.4696"
$ws.Range("E7").Value = "  +4.36%  "
$ws.Range("D8").Value = "0.3944"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "47.63"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "0.08055"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "1.024"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "21.82"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "1.871.99"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "5.941"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "7.148"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "0.00001047"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "86.69"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "0.06631"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "17.25"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "27.699.80"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "5.494"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "10.99"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "2.093.71"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "158.43"
$ws.Range("E27").Value = "  +4.36%  "
$ws.Range("D28").Value = "20.30"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").Value = "2.096"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "5.563"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "122.51"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "0.9667"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").Value = "0.09490"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").Value = "1.448"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").Value = "3.591"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "5.320"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "0.02259"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").Value = "0.06078"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "1.233"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").Value = "8.132"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "0.5976"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "0.1895"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").Value = "10.26"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.5713"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "12.18"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").Value = "3.392"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "1.939"
$ws.Range("D50").Value = "0.06853"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "114.38"
$ws.Range("E51").Value = "  +5.82%  "
